# Update countries & provincias Spain
#
# The "Pais" sheet is a flat COVID-stats table (country name in column A,
# stats in B:H) refreshed from source data. The refresh causes:
#   1) the "last updated" banner in A1 to move forward to 10:58,
#   2) a batch of per-country stat updates (new totals for several rows),
#   3) two countries (Eslovaquia, Eslovenia) to swap ranking position with
#      their neighbours, which shifts the country *label* shown on a short
#      run of rows even though their case numbers individually still line
#      up in sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) "Datos actualizados ..." banner ----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 10:58"

# ---- 2) Plain stat refreshes (country label unchanged) -------------------

# Filipinas
$ws.Range("B24").Value = 311694
$ws.Range("C24").Value = 2426
$ws.Range("D24").Value = 253488
$ws.Range("E24").Value = 52702
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 5504

# Indonesia
$ws.Range("B26").Value = 287008
$ws.Range("C26").Value = 4284
$ws.Range("D26").Value = 214947
$ws.Range("E26").Value = 61321
$ws.Range("G26").Value = 139
$ws.Range("H26").Value = 10740

# Israel
$ws.Range("B27").Value = 239806
$ws.Range("C27").Value = 2880
$ws.Range("D27").Value = 173109
$ws.Range("E27").Value = 65150
$ws.Range("G27").Value = 19
$ws.Range("H27").Value = 1547

# Barein
$ws.Range("E56").Value = 5907
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 248

# Singapur
$ws.Range("B59").Value = 57765
$ws.Range("C59").Value = 23
$ws.Range("E59").Value = 272

# Croacia
$ws.Range("B89").Value = 16593
$ws.Range("C89").Value = 213
$ws.Range("D89").Value = 15057
$ws.Range("E89").Value = 1256
$ws.Range("G89").Value = 5
$ws.Range("H89").Value = 280

# Lituania
$ws.Range("B131").Value = 4693
$ws.Range("C131").Value = 115
$ws.Range("D131").Value = 2365
$ws.Range("E131").Value = 2236

# Estonia
$ws.Range("B143").Value = 3371
$ws.Range("C143").Value = 57
$ws.Range("D143").Value = 2605
$ws.Range("E143").Value = 702

# Nueva Caledonia
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0

# ---- 3a) Jordania/Eslovaquia block reshuffles rows 103-107 ---------------
# Eslovaquia jumps ahead of Jordania, Guayana Francesa, Finlandia and
# Tayikistan, pushing each of those down one row with their own refreshed
# totals carried along.

$ws.Range("A103").Value = "Eslovaquia"
$ws.Range("B103").Value = 10141
$ws.Range("C103").Value = 567
$ws.Range("D103").Value = 4395
$ws.Range("E103").Value = 5698
$ws.Range("G103").Value = 3
$ws.Range("H103").Value = 48

$ws.Range("A104").Value = "Jordania"
$ws.Range("B104").Value = 10049
$ws.Range("D104").Value = 4496
$ws.Range("E104").Value = 5496
$ws.Range("H104").Value = 57

$ws.Range("A105").Value = "Guayana Francesa"
$ws.Range("B105").Value = 9929
$ws.Range("D105").Value = 9569
$ws.Range("E105").Value = 294
$ws.Range("H105").Value = 66

$ws.Range("A106").Value = "Finlandia"
$ws.Range("B106").Value = 9892
$ws.Range("D106").Value = 7850
$ws.Range("E106").Value = 1697
$ws.Range("H106").Value = 345

$ws.Range("A107").Value = "Tayikistan"
$ws.Range("B107").Value = 9726
$ws.Range("D107").Value = 8531
$ws.Range("E107").Value = 1120
$ws.Range("H107").Value = 75

# ---- 3b) Cuba/Eslovenia swap (rows 119-120) -------------------------------
$ws.Range("A119").Value = "Eslovenia"
$ws.Range("B119").Value = 5690
$ws.Range("C119").Value = 203
$ws.Range("D119").Value = 3804
$ws.Range("E119").Value = 1736
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 150

$ws.Range("A120").Value = "Cuba"
$ws.Range("B120").Value = 5531
$ws.Range("D120").Value = 4866
$ws.Range("E120").Value = 543
$ws.Range("H120").Value = 122
